{"js": "// Replace each two-digit multiplication expression in the document body\n// with its new value, per the commit's regenerated numbers.\nconst replacements = [\n  [\"56\u00d742=\", \"37\u00d772=\"],\n  [\"68\u00d718=\", \"92\u00d779=\"],\n  [\"23\u00d742=\", \"79\u00d760=\"],\n  [\"65\u00d755=\", \"94\u00d776=\"],\n  [\"33\u00d798=\", \"49\u00d792=\"],\n  [\"32\u00d719=\", \"63\u00d768=\"],\n  [\"81\u00d724=\", \"33\u00d732=\"],\n  [\"57\u00d763=\", \"94\u00d794=\"],\n  [\"84\u00d774=\", \"95\u00d777=\"],\n  [\"41\u00d772=\", \"33\u00d757=\"],\n  [\"33\u00d720=\", \"76\u00d787=\"],\n  [\"20\u00d767=\", \"29\u00d738=\"],\n  [\"16\u00d798=\", \"87\u00d786=\"],\n  [\"26\u00d713=\", \"90\u00d781=\"],\n  [\"47\u00d764=\", \"72\u00d772=\"],\n  [\"45\u00d716=\", \"58\u00d789=\"],\n  [\"77\u00d724=\", \"69\u00d736=\"],\n  [\"88\u00d733=\", \"22\u00d776=\"],\n  [\"32\u00d758=\", \"28\u00d740=\"],\n  [\"55\u00d761=\", \"83\u00d778=\"],\n  [\"93\u00d787=\", \"22\u00d796=\"],\n  [\"87\u00d769=\", \"67\u00d750=\"],\n  [\"25\u00d722=\", \"96\u00d784=\"],\n  [\"53\u00d726=\", \"33\u00d761=\"],\n  [\"14\u00d722=\", \"31\u00d790=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression in the document body\n# with its new value, per the commit's regenerated numbers.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"56\u00d742=\", \"37\u00d772=\"),\n    @(\"68\u00d718=\", \"92\u00d779=\"),\n    @(\"23\u00d742=\", \"79\u00d760=\"),\n    @(\"65\u00d755=\", \"94\u00d776=\"),\n    @(\"33\u00d798=\", \"49\u00d792=\"),\n    @(\"32\u00d719=\", \"63\u00d768=\"),\n    @(\"81\u00d724=\", \"33\u00d732=\"),\n    @(\"57\u00d763=\", \"94\u00d794=\"),\n    @(\"84\u00d774=\", \"95\u00d777=\"),\n    @(\"41\u00d772=\", \"33\u00d757=\"),\n    @(\"33\u00d720=\", \"76\u00d787=\"),\n    @(\"20\u00d767=\", \"29\u00d738=\"),\n    @(\"16\u00d798=\", \"87\u00d786=\"),\n    @(\"26\u00d713=\", \"90\u00d781=\"),\n    @(\"47\u00d764=\", \"72\u00d772=\"),\n    @(\"45\u00d716=\", \"58\u00d789=\"),\n    @(\"77\u00d724=\", \"69\u00d736=\"),\n    @(\"88\u00d733=\", \"22\u00d776=\"),\n    @(\"32\u00d758=\", \"28\u00d740=\"),\n    @(\"55\u00d761=\", \"83\u00d778=\"),\n    @(\"93\u00d787=\", \"22\u00d796=\"),\n    @(\"87\u00d769=\", \"67\u00d750=\"),\n    @(\"25\u00d722=\", \"96\u00d784=\"),\n    @(\"53\u00d726=\", \"33\u00d761=\"),\n    @(\"14\u00d722=\", \"31\u00d790=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute([ref]$old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$new, 2) | Out-Null\n}\n"}
